$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows (5-7) in column A, mirroring the existing rows
# (2-4), which hold "SCRIPT/.../*.ssb" filenames styled like the rest of
# column A (wrap-text, ~43.2pt tall single-line rows).
$ws.Range("A5").Value = "SCRIPT/P02P01A/um1101.ssb"
$ws.Range("A6").Value = "SCRIPT/P02P01A/um1104.ssb"
$ws.Range("A7").Value = "SCRIPT/P02P01A/um1107.ssb"

$ws.Rows(5).RowHeight = 43.2
$ws.Rows(6).RowHeight = 43.2
$ws.Rows(7).RowHeight = 43.2

# Move the active selection to D5, as in the edited workbook.
[void]$ws.Range("D5").Select()
